# The presentation ships two theme parts:
#   theme1.xml -> "Integral" theme / "Red Violet" colour scheme (used by the slide master)
#   theme2.xml -> "Office Theme" / "Office" colour scheme (used by the notes master)
# The authored edit swaps their contents, so the slide master ends up using the
# "Office Theme" palette and the notes master ends up using the "Integral" palette.
#
# The exposed PowerPoint object model only lets us edit the slide master's
# theme colour scheme (Master.Theme.ThemeColorScheme); that is the theme that
# actually renders on every slide, so we reproduce that half of the swap here,
# updating every colour slot to the target "Office Theme" palette.
# (MsoThemeColorSchemeIndex order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6,
# 11 hlink, 12 folHlink. RGB values are encoded as R + G*256 + B*65536.)

$p = $ppt.ActivePresentation
$scheme = $p.SlideMaster.Theme.ThemeColorScheme

$scheme.Item(1).RGB  = 0        + (0   * 256) + (0   * 65536)   # dk1      000000
$scheme.Item(2).RGB  = 255      + (255 * 256) + (255 * 65536)   # lt1      FFFFFF
$scheme.Item(3).RGB  = 68       + (84  * 256) + (106 * 65536)   # dk2      44546A
$scheme.Item(4).RGB  = 231      + (230 * 256) + (230 * 65536)   # lt2      E7E6E6
$scheme.Item(5).RGB  = 91       + (155 * 256) + (213 * 65536)   # accent1  5B9BD5
$scheme.Item(6).RGB  = 237      + (125 * 256) + (49  * 65536)   # accent2  ED7D31
$scheme.Item(7).RGB  = 165      + (165 * 256) + (165 * 65536)   # accent3  A5A5A5
$scheme.Item(8).RGB  = 255      + (192 * 256) + (0   * 65536)   # accent4  FFC000
$scheme.Item(9).RGB  = 68       + (114 * 256) + (196 * 65536)   # accent5  4472C4
$scheme.Item(10).RGB = 112      + (173 * 256) + (71  * 65536)   # accent6  70AD47
$scheme.Item(11).RGB = 5        + (99  * 256) + (193 * 65536)   # hlink    0563C1
$scheme.Item(12).RGB = 149      + (79  * 256) + (114 * 65536)   # folHlink 954F72
